$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("A19").Value = "SB"
$ws.Range("B19").Value = 43384
$ws.Range("C19").Value = "add a function to automate building a directory structure"

# Row 20
$ws.Range("A20").Value = "SB"
$ws.Range("B20").Value = 43384
$ws.Range("C20").Value = "take the remove_rows_cols functino from RIVAS - adds in options to just do rows or just do cols"

# Copy the existing date format (from B2, which uses numFmtId 14 / style index 1)
# onto the two new date cells so no new style entries are created.
$ws.Range("B2").Copy()
$ws.Range("B19:B20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update selection to match the new active cell after data entry
$ws.Range("A21").Select()
